# Apply the "Uploading newest EPS-US files" edit to the
# Recipient Heat Fuel Fractions (RHFF) workbook.
#
# Sheet "About": replace the EU-centric notes with new US-centric notes.
# Sheet "RHFF": the shift-to-other-fuels policy now routes everything to
#   hydrogen instead of being split across electricity/hydrogen, so the
#   "electricity" row becomes all zeros and the "hydrogen" row becomes all
#   ones.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "About" sheet
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Rows 1-7 keep their existing text; nothing to change there.

# Row 8 stays blank.

# Row 9: new sentence, and it picks up the sheet's default (unstyled) look.
$about.Range("A9").Font.Bold = $false
$about.Range("A9").Value = "District heat facilities most likely need to burn thermal fuels to generate heat (as using"

# Row 10: continuation of the sentence started on row 9.
$about.Range("A10").Value = "electricity would be inefficient, relative to the end users using electricity for heat directly)."

# Row 11 becomes blank (previously held a sentence fragment).
$about.Range("A11").Value = ""

# Row 12: new sentence.
$about.Range("A12").Value = "For the United States, where most district heat facilities currently burn natural gas, we"

# Row 13: continuation.
$about.Range("A13").Value = "specify this policy lever as a shift to hydrogen, as one of the only thermal fuel options"

# Row 14: conclusion.
$about.Range("A14").Value = "to further reduce GHG emissions."

# ---------------------------------------------------------------------
# "RHFF" sheet
# ---------------------------------------------------------------------
$rhff = $wb.Worksheets.Item("RHFF")

# The "electricity" row (row 2) no longer receives any shifted fuel use.
$rhff.Range("B2:K2").Value = 0

# A few remaining non-zero fractions in the "natural gas" (row 4) and
# "biomass" (row 5) rows are zeroed out as well.
$rhff.Cells.Item(4, 3).Value = 0   # C4 natural gas <- coal
$rhff.Cells.Item(4, 6).Value = 0   # F4 natural gas <- petroleum diesel
$rhff.Cells.Item(4, 8).Value = 0   # H4 natural gas <- crude oil
$rhff.Cells.Item(4, 9).Value = 0   # I4 natural gas <- heavy or residual fuel oil
$rhff.Cells.Item(4, 10).Value = 0  # J4 natural gas <- LPG propane or butane

$rhff.Cells.Item(5, 3).Value = 0   # C5 biomass <- coal
$rhff.Cells.Item(5, 4).Value = 0   # D5 biomass <- natural gas
$rhff.Cells.Item(5, 6).Value = 0   # F5 biomass <- petroleum diesel
$rhff.Cells.Item(5, 8).Value = 0   # H5 biomass <- crude oil
$rhff.Cells.Item(5, 9).Value = 0   # I5 biomass <- heavy or residual fuel oil
$rhff.Cells.Item(5, 10).Value = 0  # J5 biomass <- LPG propane or butane
$rhff.Cells.Item(5, 11).Value = 0  # K5 biomass <- hydrogen

# The "hydrogen" row (row 11) now receives all of the shifted fuel use.
$rhff.Range("B11:K11").Value = 1
